# Auto-generated edit script applying numeric corrections to Tiamat_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 192.25
$ws.Range("I15").Value = 192.25
$ws.Range("K15").Value = 576.75
$ws.Range("M15").Value = -407.75

$ws.Range("H132").Value = 141651.31
$ws.Range("I132").Value = 2486
$ws.Range("J132").Value = 670479.4399999999
$ws.Range("K132").Value = 7458
$ws.Range("L132").Value = 2011438.32
$ws.Range("M132").Value = -4928
$ws.Range("N132").Value = -2016498.32

$ws.Range("H137").Value = 39104.45
$ws.Range("I137").Value = 63608
$ws.Range("J137").Value = 8946.23
$ws.Range("K137").Value = 190824
$ws.Range("L137").Value = 26838.69
$ws.Range("M137").Value = -188274
$ws.Range("N137").Value = -31938.69

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1306695.9
$ws.Range("I32").Value = 1545000.9
$ws.Range("J32").Value = 15877.167
$ws.Range("K32").Value = 1545000.9
$ws.Range("L32").Value = 15877.167
$ws.Range("M32").Value = -1544713.9
$ws.Range("N32").Value = -16451.167

$ws.Range("H61").Value = 1403.1702
$ws.Range("I61").Value = 801.0294
$ws.Range("J61").Value = 2978
$ws.Range("K61").Value = 801.0294
$ws.Range("L61").Value = 2978
$ws.Range("M61").Value = -589.0294
$ws.Range("N61").Value = -3402

$ws.Range("H74").Value = 31452.059
$ws.Range("I74").Value = 40388.848
$ws.Range("J74").Value = 2407.5
$ws.Range("K74").Value = 40388.848
$ws.Range("L74").Value = 2407.5
$ws.Range("M74").Value = -39514.848
$ws.Range("N74").Value = -4155.5

$ws.Range("H77").Value = 31452.059
$ws.Range("I77").Value = 40388.848
$ws.Range("J77").Value = 2407.5
$ws.Range("K77").Value = 201944.24
$ws.Range("L77").Value = 12037.5
$ws.Range("M77").Value = -197576.24
$ws.Range("N77").Value = -20773.5

$ws.Range("H88").Value = 29986.428
$ws.Range("I88").Value = 1335.3334
$ws.Range("J88").Value = 51474.75
$ws.Range("K88").Value = 1335.3334
$ws.Range("L88").Value = 51474.75
$ws.Range("M88").Value = -929.3334
$ws.Range("N88").Value = -52286.75

$ws.Range("H91").Value = 29986.428
$ws.Range("I91").Value = 1335.3334
$ws.Range("J91").Value = 51474.75
$ws.Range("K91").Value = 1335.3334
$ws.Range("L91").Value = 51474.75
$ws.Range("M91").Value = 68.66660000000002
$ws.Range("N91").Value = -54282.75

$ws.Range("H132").Value = 1559383.5
$ws.Range("I132").Value = 1824182.1
$ws.Range("J132").Value = 632588.3
$ws.Range("K132").Value = 5472546.300000001
$ws.Range("L132").Value = 1897764.9
$ws.Range("M132").Value = -5470016.300000001
$ws.Range("N132").Value = -1902824.9

$ws.Range("H136").Value = 1403.1702
$ws.Range("I136").Value = 801.0294
$ws.Range("J136").Value = 2978
$ws.Range("K136").Value = 2403.0882
$ws.Range("L136").Value = 8934
$ws.Range("M136").Value = 146.9117999999999
$ws.Range("N136").Value = -14034

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 27873.357
$ws.Range("I134").Value = 1447.1936
$ws.Range("J134").Value = 102347.09
$ws.Range("K134").Value = 4341.5808
$ws.Range("L134").Value = 307041.27
$ws.Range("M134").Value = -1806.5808
$ws.Range("N134").Value = -312111.27

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 2166.3333
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 2999.5
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 2999.5
$ws.Range("M17").Value = -326
$ws.Range("N17").Value = -3347.5

$ws.Range("H31").Value = 10174.022
$ws.Range("I31").Value = 8788.843999999999
$ws.Range("J31").Value = 13867.833
$ws.Range("K31").Value = 8788.843999999999
$ws.Range("L31").Value = 13867.833
$ws.Range("M31").Value = -8493.843999999999
$ws.Range("N31").Value = -14457.833

$ws.Range("H34").Value = 10174.022
$ws.Range("I34").Value = 8788.843999999999
$ws.Range("J34").Value = 13867.833
$ws.Range("K34").Value = 8788.843999999999
$ws.Range("L34").Value = 13867.833
$ws.Range("M34").Value = -8586.843999999999
$ws.Range("N34").Value = -14271.833

$ws.Range("H58").Value = 1189.1538
$ws.Range("I58").Value = 752.5
$ws.Range("J58").Value = 1887.8
$ws.Range("K58").Value = 752.5
$ws.Range("L58").Value = 1887.8
$ws.Range("M58").Value = -549.5
$ws.Range("N58").Value = -2293.8

$ws.Range("H124").Value = 19649
$ws.Range("J124").Value = 19649
$ws.Range("L124").Value = 19649
$ws.Range("N124").Value = -24559

$ws.Range("H132").Value = 1364.8292
$ws.Range("I132").Value = 853.11536
$ws.Range("J132").Value = 2251.8
$ws.Range("K132").Value = 2559.34608
$ws.Range("L132").Value = 6755.400000000001
$ws.Range("M132").Value = -29.34608000000026
$ws.Range("N132").Value = -11815.4

$ws.Range("H134").Value = 1079.1321
$ws.Range("I134").Value = 1049.919
$ws.Range("K134").Value = 3149.757000000001
$ws.Range("M134").Value = -614.7570000000005

$ws.Range("H136").Value = 1189.1538
$ws.Range("I136").Value = 752.5
$ws.Range("J136").Value = 1887.8
$ws.Range("K136").Value = 2257.5
$ws.Range("L136").Value = 5663.4
$ws.Range("M136").Value = 292.5
$ws.Range("N136").Value = -10763.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 10669.3
$ws.Range("I5").Value = 707.6
$ws.Range("K5").Value = 2122.8
$ws.Range("M5").Value = -2010.8

$ws.Range("H131").Value = 56818790
$ws.Range("I131").Value = 438.8
$ws.Range("J131").Value = 104167420
$ws.Range("K131").Value = 1316.4
$ws.Range("L131").Value = 312502260
$ws.Range("M131").Value = 3723.6
$ws.Range("N131").Value = -312512340

$ws.Range("H135").Value = 10669.3
$ws.Range("I135").Value = 707.6
$ws.Range("K135").Value = 6368.400000000001
$ws.Range("M135").Value = -3833.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 27113
$ws.Range("J125").Value = 27113
$ws.Range("L125").Value = 27113
$ws.Range("N125").Value = -32033

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 418745.84
$ws.Range("I136").Value = 557180
$ws.Range("J136").Value = 3443.3333
$ws.Range("K136").Value = 1671540
$ws.Range("L136").Value = 10329.9999
$ws.Range("M136").Value = -1668990
$ws.Range("N136").Value = -15429.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2561.4575
$ws.Range("I132").Value = 619.02325
$ws.Range("J132").Value = 7781.75
$ws.Range("K132").Value = 2559.34608
$ws.Range("L132").Value = 23345.25
$ws.Range("M132").Value = 672.9302500000001
$ws.Range("N132").Value = -28405.25
